$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3652
$ws.Range("C2").Value = 0.4999999999999999
$ws.Range("D2").Value = 0.1939
$ws.Range("E2").Value = 0.3794000000000001
$ws.Range("F2").Value = 0.56296
$ws.Range("G2").Value = 0.2232
$ws.Range("H2").Value = 0.2518

$ws.Range("B3").Value = 0.3452
$ws.Range("C3").Value = 0.4999999999999999
$ws.Range("D3").Value = 0.2979
$ws.Range("E3").Value = 0.3794000000000001
$ws.Range("F3").Value = 0.45896
$ws.Range("G3").Value = 0.2432
$ws.Range("H3").Value = 0.2518

$ws.Range("B4").Value = 0.3452
$ws.Range("C4").Value = 0.4999999999999999
$ws.Range("D4").Value = 0.2979
$ws.Range("E4").Value = 0.3794000000000001
$ws.Range("F4").Value = 0.4789600000000001
$ws.Range("G4").Value = 0.2232
$ws.Range("H4").Value = 0.2518
